# Applies the commit's changes to the "Blank 3 Statement Model" sheet:
#   1. Reword the Balance Sheet Check label in A3.
#   2. Hide column E and set its width to match the template's narrow
#      "spacer/helper" column styling (width ~13 chars).
#   3. Replace the Retained Earnings roll-forward formulas in C77/D77
#      (and the B77 seed value) with flat zero values, matching the
#      "blank template" reset baked into this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blank 3 Statement Model")

# 1) Update the Balance Sheet Check label.
$ws.Range("A3").Value = "Balance Sheet Check (A - (L + E))"

# 2) Hide column E with width 13 (raw OOXML width == ColumnWidth + 5/6).
$ws.Columns.Item(5).ColumnWidth = 12.166666666666666
$ws.Columns.Item(5).Hidden = $true

# 3) Reset Retained Earnings row (B77:D77) to static zeros, recolouring
#    C77/D77's font to the same blue "hardcoded input" colour already used
#    by B77 (style 105 vs. the old formula style 97 -- they differ only by
#    font colour; the number format is already shared).
$ws.Range("C77").Font.Color = $ws.Range("B77").Font.Color
$ws.Range("D77").Font.Color = $ws.Range("B77").Font.Color

$ws.Range("B77").Value = 0
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 0
